{"js": "// Replace the three-digit x one-digit multiplication problems with new ones,\n// per the commit diff. Each old equation text is unique in the document, so\n// a simple exact-match search/replace for each pair is unambiguous.\nconst replacements = [\n  [\"593\u00d75=2965\", \"170\u00d79=1530\"],\n  [\"817\u00d72=1634\", \"864\u00d78=6912\"],\n  [\"232\u00d78=1856\", \"336\u00d72=672\"],\n  [\"133\u00d73=399\", \"530\u00d74=2120\"],\n  [\"146\u00d78=1168\", \"667\u00d73=2001\"],\n  [\"964\u00d75=4820\", \"185\u00d72=370\"],\n  [\"879\u00d79=7911\", \"436\u00d78=3488\"],\n  [\"739\u00d73=2217\", \"467\u00d76=2802\"],\n  [\"436\u00d76=2616\", \"145\u00d79=1305\"],\n  [\"485\u00d79=4365\", \"320\u00d79=2880\"],\n  [\"171\u00d76=1026\", \"177\u00d77=1239\"],\n  [\"961\u00d77=6727\", \"279\u00d72=558\"],\n  [\"443\u00d77=3101\", \"312\u00d78=2496\"],\n  [\"715\u00d73=2145\", \"849\u00d76=5094\"],\n  [\"826\u00d74=3304\", \"556\u00d72=1112\"],\n  [\"876\u00d75=4380\", \"699\u00d75=3495\"],\n  [\"119\u00d73=357\", \"384\u00d75=1920\"],\n  [\"553\u00d78=4424\", \"987\u00d75=4935\"],\n  [\"128\u00d76=768\", \"759\u00d72=1518\"],\n  [\"800\u00d79=7200\", \"599\u00d75=2995\"],\n  [\"288\u00d78=2304\", \"783\u00d78=6264\"],\n  [\"957\u00d75=4785\", \"421\u00d75=2105\"],\n  [\"465\u00d74=1860\", \"503\u00d74=2012\"],\n  [\"490\u00d72=980\", \"712\u00d75=3560\"],\n  [\"591\u00d72=1182\", \"149\u00d76=894\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems with new ones,\n# per the commit diff. Each old equation text is unique in the document, so\n# an exact-match Find/Replace for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"593\u00d75=2965\", \"170\u00d79=1530\"),\n    @(\"817\u00d72=1634\", \"864\u00d78=6912\"),\n    @(\"232\u00d78=1856\", \"336\u00d72=672\"),\n    @(\"133\u00d73=399\", \"530\u00d74=2120\"),\n    @(\"146\u00d78=1168\", \"667\u00d73=2001\"),\n    @(\"964\u00d75=4820\", \"185\u00d72=370\"),\n    @(\"879\u00d79=7911\", \"436\u00d78=3488\"),\n    @(\"739\u00d73=2217\", \"467\u00d76=2802\"),\n    @(\"436\u00d76=2616\", \"145\u00d79=1305\"),\n    @(\"485\u00d79=4365\", \"320\u00d79=2880\"),\n    @(\"171\u00d76=1026\", \"177\u00d77=1239\"),\n    @(\"961\u00d77=6727\", \"279\u00d72=558\"),\n    @(\"443\u00d77=3101\", \"312\u00d78=2496\"),\n    @(\"715\u00d73=2145\", \"849\u00d76=5094\"),\n    @(\"826\u00d74=3304\", \"556\u00d72=1112\"),\n    @(\"876\u00d75=4380\", \"699\u00d75=3495\"),\n    @(\"119\u00d73=357\", \"384\u00d75=1920\"),\n    @(\"553\u00d78=4424\", \"987\u00d75=4935\"),\n    @(\"128\u00d76=768\", \"759\u00d72=1518\"),\n    @(\"800\u00d79=7200\", \"599\u00d75=2995\"),\n    @(\"288\u00d78=2304\", \"783\u00d78=6264\"),\n    @(\"957\u00d75=4785\", \"421\u00d75=2105\"),\n    @(\"465\u00d74=1860\", \"503\u00d74=2012\"),\n    @(\"490\u00d72=980\", \"712\u00d75=3560\"),\n    @(\"591\u00d72=1182\", \"149\u00d76=894\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
